# Update "想去人数" (want-to-go count) values in column F across sheets,
# reflecting the data refresh captured in the gh-pages output at 456a3b4.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibition) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 1119
$ws1.Range("F4").Value = 1810
$ws1.Range("F5").Value = 795
$ws1.Range("F6").Value = 403
$ws1.Range("F7").Value = 218

# --- Sheet "演出" (Performance) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 14

# --- Sheet "全部类型" (All types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 1119
$ws4.Range("F4").Value = 1810
$ws4.Range("F5").Value = 14
$ws4.Range("F6").Value = 795
$ws4.Range("F7").Value = 403
$ws4.Range("F8").Value = 218
